# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the style from an existing header cell (A1) onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-56): same record (68-94-0) repeated for every player ---
for ($r = 2; $r -le 56; $r++) {
    $ws.Range("AD$r").Value = 68
    $ws.Range("AE$r").Value = 94
    $ws.Range("AF$r").Value = 0
}
